$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  34"
$ws.Range("C9").Value = "Report Covering the Week  8/21/2023  Through  8/27/2023"

# --- Cells that change data type (string <-> number) need style fix-up via PasteSpecial ---
$ws.Range("C15").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = 4

$ws.Range("E15").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = -50

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("C15").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1

$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F28").PasteSpecial(-4122)

$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F29").PasteSpecial(-4122)

$ws.Range("C15").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = 1

$ws.Range("E15").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100

$ws.Range("C15").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("G30").Value = 1

$ws.Range("E15").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("H30").Value = -100

$excel.CutCopyMode = $false

# --- Plain numeric value updates ---
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 7
$ws.Range("J15").Value = 11
$ws.Range("K15").Value = -36.363636363636
$ws.Range("L15").Value = 75
$ws.Range("M15").Value = -30
$ws.Range("N15").Value = -61.111111111111
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 12
$ws.Range("H16").Value = -20
$ws.Range("I16").Value = 94
$ws.Range("J16").Value = 73
$ws.Range("K16").Value = 28.767123287671
$ws.Range("L16").Value = 123.809523809524
$ws.Range("M16").Value = -19.658119658119
$ws.Range("N16").Value = -81.640625
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 27.777777777777
$ws.Range("I17").Value = 122
$ws.Range("J17").Value = 106
$ws.Range("K17").Value = 15.094339622641
$ws.Range("L17").Value = 25.773195876288
$ws.Range("M17").Value = 46.987951807228
$ws.Range("N17").Value = -37.435897435897
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -25
$ws.Range("I18").Value = 62
$ws.Range("J18").Value = 81
$ws.Range("K18").Value = -23.456790123456
$ws.Range("L18").Value = 63.157894736842
$ws.Range("M18").Value = 21.56862745098
$ws.Range("N18").Value = -84.803921568627
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 16.666666666666
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 13.333333333333
$ws.Range("I19").Value = 261
$ws.Range("J19").Value = 234
$ws.Range("K19").Value = 11.538461538461
$ws.Range("L19").Value = 76.351351351351
$ws.Range("M19").Value = 46.629213483146
$ws.Range("N19").Value = -39.443155452436
$ws.Range("C20").Value = 2
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -20
$ws.Range("I20").Value = 74
$ws.Range("J20").Value = 48
$ws.Range("K20").Value = 54.166666666666
$ws.Range("L20").Value = 146.666666666667
$ws.Range("M20").Value = 362.5
$ws.Range("N20").Value = -70.980392156862
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = 5.263157894736
$ws.Range("F21").Value = 85
$ws.Range("G21").Value = 82
$ws.Range("H21").Value = 3.658536585365
$ws.Range("I21").Value = 620
$ws.Range("J21").Value = 554
$ws.Range("K21").Value = 11.913357400722
$ws.Range("L21").Value = 72.222222222222
$ws.Range("M21").Value = 35.667396061269
$ws.Range("N21").Value = -66.064586754241
$ws.Range("C22").Value = 2
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 400
$ws.Range("I22").Value = 26
$ws.Range("K22").Value = 100
$ws.Range("L22").Value = 52.941176470588
$ws.Range("M22").Value = 23.809523809523
$ws.Range("C23").Value = 4
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 21
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 133.333333333333
$ws.Range("I23").Value = 104
$ws.Range("J23").Value = 80
$ws.Range("K23").Value = 30
$ws.Range("L23").Value = 44.444444444444
$ws.Range("M23").Value = 76.271186440678
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 8
$ws.Range("E24").Value = -12.5
$ws.Range("F24").Value = 32
$ws.Range("G24").Value = 36
$ws.Range("H24").Value = -11.111111111111
$ws.Range("I24").Value = 330
$ws.Range("J24").Value = 331
$ws.Range("K24").Value = -0.302114803625
$ws.Range("L24").Value = 25.475285171102
$ws.Range("M24").Value = -19.117647058823
$ws.Range("C25").Value = 4
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = 7.142857142857
$ws.Range("I25").Value = 165
$ws.Range("J25").Value = 150
$ws.Range("K25").Value = 10
$ws.Range("L25").Value = 32
$ws.Range("M25").Value = -17.910447761194
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 8
$ws.Range("J26").Value = 17
$ws.Range("K26").Value = -52.941176470588
$ws.Range("L26").Value = 60
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 28
$ws.Range("J27").Value = 33
$ws.Range("K27").Value = -15.151515151515
$ws.Range("L27").Value = 21.739130434782
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = -100
$ws.Range("N28").Value = -68.181818181818
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = -100
$ws.Range("N29").Value = -61.111111111111
$ws.Range("J30").Value = 8
$ws.Range("K30").Value = -75
